$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 1 to make room for the new texture/icon rows
$insertRange = $ws.Range("A1:A3").EntireRow
$insertRange.Insert()

# Row 1: texture / png / default.png
$ws.Range("B1").Value = "png"
$ws.Range("A1").Value = "texture"
$ws.Range("C1").Value = "default.png"

# Row 2: texture / png / road.png
$ws.Range("A2").Value = "texture"
$ws.Range("B2").Value = "png"
$ws.Range("C2").Value = "road.png"

# Row 3: icon / png / icon.png
$ws.Range("A3").Value = "icon"
$ws.Range("B3").Value = "png"
$ws.Range("C3").Value = "icon.png"

# Update selection to C3
$ws.Range("C3").Select()
